$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '304.73'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.37%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '10'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.92'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.47%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '10'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.095'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.96%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '10'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08082'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.76%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '10'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.922'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.01%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '10'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.743'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.22%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '10'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9276'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.40%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '10'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1346'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.14%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '10'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1906'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.56%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '10'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09178'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-4.57%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '10'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03402'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-5.80%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '10'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09826'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.33%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '10'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001409'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.63%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '10'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005905'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.61%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '10'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.555'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.67%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '10'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.180'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '3.50%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '10'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.015'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.18%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '10'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3452'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.67%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '10'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.63%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '10'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.906'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.85%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '10'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '5.24%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '10'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04440'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.97%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '10'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.13%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '10'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004811'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.28%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '10'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001301'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '3.79%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '10'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '4.14%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '10'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '10'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '10'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '10'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '10'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '10'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '10'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '10'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '10'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '10'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '10'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '10'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02007'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '6.24%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '10'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04908'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4.85%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '10'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007646'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.23%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '10'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009974'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.03%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '10'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1377'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.99%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '10'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002105'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.53%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '10'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01090'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.05%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '10'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006397'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '3.16%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '10'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.20%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '10'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '10'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '10'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.20%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '10'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.20%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '10'
